$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells G1 / H1, copying the header format (border/bold/center) from F1
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Updated MSE (B) / R2 (C) / MAE (D) values, plus new Elapsed Time (G) / CPU (H) columns for rows 2-10
$values = @(
    @(1.330325824511135, 0.7832761913694432, 0.9555549925342113),
    @(6.848323667209031, 0.9021662183526097, 1.970736775091027),
    @(3.480567636655565, 0.8280740622608915, 1.510404694356846),
    @(2.955423048646301, 0.9980587558402471, 1.265663291807137),
    @(2.112808609028797, 0.9769773194256061, 1.154332099098965),
    @(1.816420092496004, 0.9989815900158616, 1.066839266934659),
    @(2.228590673907531, 0.9974472535085364, 1.23873864560599),
    @(14.94167941894868, 0.8213825473888408, 3.184123592369467),
    @(1.546085115677007, 0.9953911681576642, 1.003074388147043)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i][0]
    $ws.Cells.Item($row, 3).Value = $values[$i][1]
    $ws.Cells.Item($row, 4).Value = $values[$i][2]
    $ws.Cells.Item($row, 7).Value = 1.127317944850074
    $ws.Cells.Item($row, 8).Value = 0.985
}
